# Updated cryptos list (price + 1h volume change columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') values that look like plain decimal numbers need a leading
# apostrophe so Excel keeps them as literal text (matching the source format,
# e.g. '36.60' rather than being normalised to the number 36.6).

$ws.Range("D2").Value = "71.070.48"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.853.73"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'695.83"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").Value = "'173.01"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").Value = "3.852.76"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  +6.31%  "
$ws.Range("D14").Value = "'36.60"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("D15").Value = "4.504.28"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "3.858.94"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "71.175.06"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'7.27"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "'17.75"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D21").Value = "'11.14"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("D22").Value = "'496.06"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'84.88"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").Value = "'12.33"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'10.57"
$ws.Range("E27").Value = "  +2.84%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "4.012.81"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("E30").Value = "  +11.34%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "'7.61"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "'29.76"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'9.30"
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("D37").Value = "3.806.30"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("E40").Value = "  +12.76%  "
$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").Value = "'6.06"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("E43").Value = "  +5.90%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D46").Value = "'164.77"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("E47").Value = "  +5.41%  "
$ws.Range("D48").Value = "'48.69"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").Value = "'44.41"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").Value = "'418.76"
$ws.Range("E50").Value = "  +7.26%  "
$ws.Range("D51").Value = "'0.302"
$ws.Range("E51").Value = "  +0.83%  "
